$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.561.33'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '3.141.95'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'237.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.35%  '
$ws.Range('D6').Value = "'644.49"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.50%  '
$ws.Range('E7').Value = '  +11.30%  '
$ws.Range('E8').Value = '  -5.19%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '3.136.86'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('D11').Value = "'0.722"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('E12').Value = '  +4.14%  '
$ws.Range('D13').Value = "'36.75"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.85%  '
$ws.Range('E14').Value = '  -4.37%  '
$ws.Range('D15').Value = "'5.65"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.80%  '
$ws.Range('D16').Value = '90.302.08'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').Value = '3.715.02'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = '3.117.58'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = "'14.52"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.30%  '
$ws.Range('D21').Value = "'0.0000216"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('D22').Value = "'450.96"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('D23').Value = "'5.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.01%  '
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range('D25').Value = "'6.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.20%  '
$ws.Range('D26').Value = "'91.42"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.67%  '
$ws.Range('E27').Value = '  +2.56%  '
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = "'10.01"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.95%  '
$ws.Range('E31').Value = '  -3.79%  '
$ws.Range('D32').Value = "'27.35"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +15.57%  '
$ws.Range('D33').Value = "'0.202"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +33.21%  '
$ws.Range('E34').Value = '  +3.98%  '
$ws.Range('D35').Value = "'518.93"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  +5.72%  '
$ws.Range('E37').Value = '  +5.49%  '
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('E39').Value = '  +2.89%  '
$ws.Range('D40').Value = "'0.424"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.40%  '
$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D41').Value = "'22.20"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = "'0.0862"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D44').Value = "'0.743"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -16.35%  '
$ws.Range('D45').Value = "'3.32"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +39.70%  '
$ws.Range('E46').Value = '  +2.18%  '
$ws.Range('E47').Value = '  +14.00%  '
$ws.Range('D48').Value = "'149.97"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.05%  '
$ws.Range('D49').Value = "'4.60"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.99%  '
$ws.Range('D50').Value = "'45.60"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('E51').Value = '  +4.67%  '
